$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1721.5
$ws.Range("I41").Value = 1514.4286
$ws.Range("K41").Value = 1514.4286
$ws.Range("M41").Value = -1074.4286

$ws.Range("H111").Value = 4491.1113
$ws.Range("I111").Value = 4198
$ws.Range("K111").Value = 12594
$ws.Range("M111").Value = -9527

$ws.Range("H132").Value = 4732.769
$ws.Range("I132").Value = 4621.5938
$ws.Range("K132").Value = 13864.7814
$ws.Range("M132").Value = -11334.7814

$ws.Range("H135").Value = 1769.9474
$ws.Range("I135").Value = 1358.125
$ws.Range("K135").Value = 12223.125
$ws.Range("M135").Value = -9688.125

$ws.Range("H137").Value = 90423.3
$ws.Range("I137").Value = 90423.3
$ws.Range("K137").Value = 271269.9
$ws.Range("M137").Value = -268719.9

$ws.Range("H138").Value = 3541.197
$ws.Range("I138").Value = 1130.9
$ws.Range("K138").Value = 3392.7
$ws.Range("M138").Value = 1747.3

$ws.Range("H141").Value = 10866.8125
$ws.Range("I141").Value = 5267.68
$ws.Range("J141").Value = 30863.715
$ws.Range("K141").Value = 15803.04
$ws.Range("L141").Value = 92591.145
$ws.Range("M141").Value = -10623.04
$ws.Range("N141").Value = -102951.145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 168076.5
$ws.Range("I2").Value = 168076.5
$ws.Range("K2").Value = 168076.5
$ws.Range("M2").Value = -167963.5

$ws.Range("H5").Value = 81.25
$ws.Range("I5").Value = 83.333336
$ws.Range("J5").Value = 75
$ws.Range("K5").Value = 83.333336
$ws.Range("L5").Value = 75
$ws.Range("M5").Value = 28.666664
$ws.Range("N5").Value = -299

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H74").Value = 248447.55
$ws.Range("I74").Value = 181579.8
$ws.Range("K74").Value = 181579.8
$ws.Range("M74").Value = -180705.8

$ws.Range("H77").Value = 248447.55
$ws.Range("I77").Value = 181579.8
$ws.Range("K77").Value = 907899
$ws.Range("M77").Value = -903531

$ws.Range("H116").Value = 168076.5
$ws.Range("I116").Value = 168076.5
$ws.Range("K116").Value = 168076.5
$ws.Range("M116").Value = -165782.5

$ws.Range("H122").Value = 3518.3333
$ws.Range("I122").Value = 3474.5454
$ws.Range("K122").Value = 10423.6362
$ws.Range("M122").Value = -7973.636200000001

$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 168076.5
$ws.Range("I3").Value = 168076.5
$ws.Range("K3").Value = 168076.5
$ws.Range("M3").Value = -167962.5

$ws.Range("H4").Value = 81.25
$ws.Range("I4").Value = 83.333336
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 83.333336
$ws.Range("L4").Value = 75
$ws.Range("M4").Value = 31.666664
$ws.Range("N4").Value = -305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1911.9524
$ws.Range("I16").Value = 1427.3846
$ws.Range("J16").Value = 2699.375
$ws.Range("K16").Value = 1427.3846
$ws.Range("L16").Value = 2699.375
$ws.Range("M16").Value = -1140.3846
$ws.Range("N16").Value = -3273.375

$ws.Range("H31").Value = 21292.36
$ws.Range("J31").Value = 33973.535
$ws.Range("L31").Value = 33973.535
$ws.Range("N31").Value = -34563.535

$ws.Range("H34").Value = 21292.36
$ws.Range("J34").Value = 33973.535
$ws.Range("L34").Value = 33973.535
$ws.Range("N34").Value = -34377.535

$ws.Range("H58").Value = 2475.2222
$ws.Range("I58").Value = 1863.1052
$ws.Range("J58").Value = 3929
$ws.Range("K58").Value = 1863.1052
$ws.Range("L58").Value = 3929
$ws.Range("M58").Value = -1660.1052
$ws.Range("N58").Value = -4335

$ws.Range("H113").Value = 1911.9524
$ws.Range("I113").Value = 1427.3846
$ws.Range("J113").Value = 2699.375
$ws.Range("K113").Value = 1427.3846
$ws.Range("L113").Value = 2699.375
$ws.Range("M113").Value = 742.6153999999999
$ws.Range("N113").Value = -7039.375

$ws.Range("H134").Value = 28867.027
$ws.Range("I134").Value = 44462.09
$ws.Range("J134").Value = 4360.5
$ws.Range("K134").Value = 133386.27
$ws.Range("L134").Value = 13081.5
$ws.Range("M134").Value = -130851.27
$ws.Range("N134").Value = -18151.5

$ws.Range("H136").Value = 2475.2222
$ws.Range("I136").Value = 1863.1052
$ws.Range("J136").Value = 3929
$ws.Range("K136").Value = 5589.3156
$ws.Range("L136").Value = 11787
$ws.Range("M136").Value = -3039.3156
$ws.Range("N136").Value = -16887

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4076.1904
$ws.Range("J113").Value = 2064.7058
$ws.Range("L113").Value = 6194.117400000001
$ws.Range("N113").Value = -10534.1174

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 16000
$ws.Range("I35").Value = 16000
$ws.Range("K35").Value = 16000
$ws.Range("M35").Value = -15702

$ws.Range("H102").Value = 32706.576
$ws.Range("I102").Value = 1689.25
$ws.Range("J102").Value = 115419.445
$ws.Range("K102").Value = 1689.25
$ws.Range("L102").Value = 115419.445
$ws.Range("M102").Value = -67.25
$ws.Range("N102").Value = -118663.445

$ws.Range("H132").Value = 3304.4814
$ws.Range("I132").Value = 3239.6538
$ws.Range("J132").Value = 4990
$ws.Range("K132").Value = 9718.9614
$ws.Range("L132").Value = 14970
$ws.Range("M132").Value = -7188.9614
$ws.Range("N132").Value = -20030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7435
$ws.Range("I61").Value = 7401.7646
$ws.Range("K61").Value = 7401.7646
$ws.Range("M61").Value = -7199.7646

$ws.Range("H93").Value = 1915.909
$ws.Range("I93").Value = 1612.4286
$ws.Range("K93").Value = 1612.4286
$ws.Range("M93").Value = -364.4286

$ws.Range("H113").Value = 7435
$ws.Range("I113").Value = 7401.7646
$ws.Range("K113").Value = 7401.7646
$ws.Range("M113").Value = -5231.7646

$ws.Range("H136").Value = 29570.275
$ws.Range("I136").Value = 45019.418
$ws.Range("K136").Value = 135058.254
$ws.Range("M136").Value = -132508.254

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4030.3333
$ws.Range("I81").Value = 4030.3333
$ws.Range("K81").Value = 8060.6666
$ws.Range("M81").Value = -6999.6666

$ws.Range("H84").Value = 4030.3333
$ws.Range("I84").Value = 4030.3333
$ws.Range("K84").Value = 40303.333
$ws.Range("M84").Value = -34999.333

$ws.Range("H107").Value = 1360.4
$ws.Range("I107").Value = 1200.5
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 3601.5
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -1681.5
$ws.Range("N107").Value = -9840

$ws.Range("H132").Value = 19245.928
$ws.Range("I132").Value = 3625.3555
$ws.Range("K132").Value = 10876.0665
$ws.Range("M132").Value = -8346.066500000001
